$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46061
$ws.Range("D2").Value = 106
$ws.Range("E2").Value = 29
$ws.Range("G2").Value = 29
$ws.Range("M2").Value = 52
$ws.Range("N2").Value = 27
$ws.Range("O2").Value = 109

$ws.Range("A3").Value = 46061
$ws.Range("D3").Value = 145
$ws.Range("E3").Value = 39
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 39
$ws.Range("H3").Value = 4
$ws.Range("M3").Value = 59
$ws.Range("N3").Value = 63
$ws.Range("O3").Value = 162

$ws.Range("A4").Value = 46061
$ws.Range("D4").Value = 44
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 2
$ws.Range("M4").Value = 67
$ws.Range("N4").Value = 26
$ws.Range("O4").Value = 120

$ws.Range("A5").Value = 46061
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 41
$ws.Range("G5").Value = 41
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 30
$ws.Range("O5").Value = 137

$ws.Range("A6").Value = 46061
$ws.Range("E6").Value = 24
$ws.Range("G6").Value = 24
$ws.Range("M6").Value = 77
$ws.Range("N6").Value = 13
$ws.Range("O6").Value = 115

$ws.Range("A7").Value = 46061
$ws.Range("D7").Value = 173
$ws.Range("E7").Value = 26
$ws.Range("G7").Value = 26
$ws.Range("M7").Value = 44
$ws.Range("N7").Value = 41
$ws.Range("O7").Value = 112

$ws.Range("A8").Value = 46061
$ws.Range("D8").Value = 192
$ws.Range("F8").Value = 3
$ws.Range("H8").Value = 3
$ws.Range("M8").Value = 48
$ws.Range("N8").Value = 53
$ws.Range("O8").Value = 136

$ws.Range("A9").Value = 46061
$ws.Range("D9").Value = 69
$ws.Range("E9").Value = 34
$ws.Range("G9").Value = 34
$ws.Range("M9").Value = 74
$ws.Range("N9").Value = 41
$ws.Range("O9").Value = 150

$ws.Range("A10").Value = 46061
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 26
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 26
$ws.Range("H10").Value = 1
$ws.Range("M10").Value = 62
$ws.Range("N10").Value = 10
$ws.Range("O10").Value = 99

$ws.Range("A11").Value = 46061
$ws.Range("D11").Value = 134
$ws.Range("E11").Value = 34
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 34
$ws.Range("H11").Value = 5
$ws.Range("M11").Value = 59
$ws.Range("N11").Value = 43
$ws.Range("O11").Value = 135

$ws.Range("A12").Value = 46061
$ws.Range("D12").Value = 148
$ws.Range("E12").Value = 44
$ws.Range("F12").Value = 12
$ws.Range("G12").Value = 44
$ws.Range("H12").Value = 12
$ws.Range("M12").Value = 66
$ws.Range("N12").Value = 72
$ws.Range("O12").Value = 180

$ws.Range("A13").Value = 46061
$ws.Range("D13").Value = 77
$ws.Range("E13").Value = 77
$ws.Range("F13").Value = 21
$ws.Range("G13").Value = 77
$ws.Range("H13").Value = 21
$ws.Range("M13").Value = 106
$ws.Range("N13").Value = 69
$ws.Range("O13").Value = 252

$ws.Range("A14").Value = 46061
$ws.Range("D14").Value = 12
$ws.Range("E14").Value = 71
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 71
$ws.Range("H14").Value = 5
$ws.Range("M14").Value = 122
$ws.Range("N14").Value = 10
$ws.Range("O14").Value = 199

$ws.Range("A15").Value = 46061
$ws.Range("D15").Value = 115
$ws.Range("E15").Value = 23
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 23
$ws.Range("H15").Value = 2
$ws.Range("N15").Value = 58
$ws.Range("O15").Value = 139

$ws.Range("A16").Value = 46061
$ws.Range("D16").Value = 133
$ws.Range("E16").Value = 39
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 39
$ws.Range("H16").Value = 8
$ws.Range("M16").Value = 69
$ws.Range("N16").Value = 125
$ws.Range("O16").Value = 234

$ws.Range("A17").Value = 46061
$ws.Range("D17").Value = 40
$ws.Range("E17").Value = 39
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = 10
$ws.Range("M17").Value = 84
$ws.Range("N17").Value = 27
$ws.Range("O17").Value = 150

$ws.Range("A18").Value = 46061
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 44
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 44
$ws.Range("H18").Value = 6
$ws.Range("M18").Value = 75
$ws.Range("N18").Value = 13
$ws.Range("O18").Value = 132

$ws.Range("A19").Value = 46061
$ws.Range("D19").Value = 86
$ws.Range("E19").Value = 34
$ws.Range("G19").Value = 34
$ws.Range("M19").Value = 85
$ws.Range("N19").Value = 33
$ws.Range("O19").Value = 121

$ws.Range("A20").Value = 46061
$ws.Range("D20").Value = 98
$ws.Range("E20").Value = 68
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 68
$ws.Range("H20").Value = 5
$ws.Range("M20").Value = 127
$ws.Range("N20").Value = 196
$ws.Range("O20").Value = 275

$ws.Range("A21").Value = 46061
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 54
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 54
$ws.Range("H21").Value = 4
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 45
$ws.Range("O21").Value = 161

$ws.Range("A22").Value = 46061
$ws.Range("D22").Value = 14
$ws.Range("E22").Value = 50
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 50
$ws.Range("H22").Value = 5
$ws.Range("M22").Value = 78
$ws.Range("N22").Value = 33
$ws.Range("O22").Value = 127

$ws.Range("A23").Value = 46061
$ws.Range("D23").Value = 123
$ws.Range("E23").Value = 33
$ws.Range("F23").Value = 2
$ws.Range("I23").Value = 33
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = 1
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 74
$ws.Range("N23").Value = 97
$ws.Range("O23").Value = 203

$ws.Range("A24").Value = 46061
$ws.Range("D24").Value = 119
$ws.Range("E24").Value = 33
$ws.Range("F24").Value = 4
$ws.Range("I24").Value = 33
$ws.Range("J24").Value = 4
$ws.Range("M24").Value = 74
$ws.Range("N24").Value = 96
$ws.Range("O24").Value = 203

$ws.Range("A25").Value = 46061
$ws.Range("D25").Value = 43
$ws.Range("E25").Value = 44
$ws.Range("F25").Value = 4
$ws.Range("I25").Value = 44
$ws.Range("J25").Value = 4
$ws.Range("M25").Value = 73
$ws.Range("N25").Value = 28
$ws.Range("O25").Value = 144

$ws.Range("A26").Value = 46061
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 31
$ws.Range("F26").Value = 2
$ws.Range("I26").Value = 31
$ws.Range("J26").Value = 2
$ws.Range("M26").Value = 79
$ws.Range("N26").Value = 10
$ws.Range("O26").Value = 107

$ws.Range("A27").Value = 46061
$ws.Range("D27").Value = 137
$ws.Range("E27").Value = 28
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 28
$ws.Range("H27").Value = 1
$ws.Range("M27").Value = 46
$ws.Range("N27").Value = 41
$ws.Range("O27").Value = 122

$ws.Range("A28").Value = 46061
$ws.Range("D28").Value = 143
$ws.Range("E28").Value = 33
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 33
$ws.Range("H28").Value = 2
$ws.Range("M28").Value = 57
$ws.Range("N28").Value = 52
$ws.Range("O28").Value = 154

$ws.Range("A29").Value = 46061
$ws.Range("D29").Value = 60
$ws.Range("E29").Value = 38
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 38
$ws.Range("H29").Value = 3
$ws.Range("M29").Value = 73
$ws.Range("N29").Value = 28
$ws.Range("O29").Value = 142

$ws.Range("A30").Value = 46061
$ws.Range("D30").Value = 9
$ws.Range("E30").Value = 27
$ws.Range("G30").Value = 27
$ws.Range("M30").Value = 71
$ws.Range("N30").Value = 17
$ws.Range("O30").Value = 111
